$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("CW1").Value = "Race Unknown"
$ws.Range("CV1").Value = "Race Other"
$ws.Range("CX1").Value = "Race Refused to Answer"
